$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B2").Value = 0.6226591760299626
$ws.Range("C2").Value = 0.5746864310148233
$ws.Range("D2").Value = 0.9438202247191011
$ws.Range("E2").Value = 0.7143869596031184
$ws.Range("F2").Value = 0.8363757052771325
$ws.Range("G2").Value = 0.921065579531876
$ws.Range("H2").Value = 0.7968410273674761
$ws.Range("I2").Value = 504
$ws.Range("J2").Value = 373
$ws.Range("K2").Value = 161
$ws.Range("L2").Value = 30

# --- Sheet: Classification Report ---
$ws = $wb.Worksheets.Item("Classification Report")
# Row 2 - class "0"
$ws.Range("B2").Value = 0.8429319371727748
$ws.Range("C2").Value = 0.301498127340824
$ws.Range("D2").Value = 0.4441379310344827

# Row 3 - class "1"
$ws.Range("B3").Value = 0.5746864310148233
$ws.Range("C3").Value = 0.9438202247191011
$ws.Range("D3").Value = 0.7143869596031184

# Row 4 - accuracy
$ws.Range("B4").Value = 0.6226591760299626
$ws.Range("C4").Value = 0.6226591760299626
$ws.Range("D4").Value = 0.6226591760299626
$ws.Range("E4").Value = 0.6226591760299626

# Row 5 - macro avg
$ws.Range("B5").Value = 0.7088091840937991
$ws.Range("C5").Value = 0.6226591760299626
$ws.Range("D5").Value = 0.5792624453188006

# Row 6 - weighted avg
$ws.Range("B6").Value = 0.7088091840937991
$ws.Range("C6").Value = 0.6226591760299626
$ws.Range("D6").Value = 0.5792624453188006

# --- Sheet: Confusion Matrix ---
$ws = $wb.Worksheets.Item("Confusion Matrix")
# Row 2 - Actual 0
$ws.Range("B2").Value = 161
$ws.Range("C2").Value = 373

# Row 3 - Actual 1
$ws.Range("B3").Value = 30
$ws.Range("C3").Value = 504
